# Apply inventory corrections to the "Inventario Quintales" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario Quintales")

# Fix typo in product name: "Ponedora fas 2" -> "Ponedora fase 2"
$ws.Range("A12").Value = "Ponedora fase 2"

# Update "Existencias" (stock) column E values
$ws.Range("E2").Value = 100
$ws.Range("E3").Value = 500
$ws.Range("E4").Value = 500
$ws.Range("E5").Value = 500
$ws.Range("E6").Value = 493
$ws.Range("E7").Value = 500
$ws.Range("E8").Value = 500
$ws.Range("E9").Value = 500
$ws.Range("E10").Value = 500
$ws.Range("E11").Value = 500
$ws.Range("E12").Value = 500
$ws.Range("E13").Value = 500
$ws.Range("E14").Value = 1300
